# Regenerate merged AHB files
# Re-applies the "group header" gray highlight style (already used on row 9)
# to the next batch of group-header rows, and clears the now-stale
# "AENDERUNG" (Change) marker out of column L for every data row in those
# groups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Rows whose entire A:V span is restyled to the gray "group header" look
# (column B additionally becomes bold). Row 9 already carries this exact
# formatting in the workbook, so it is used as the style donor.
$fullRestyleRows = @(13, 17, 23, 27, 34, 40, 71, 99)

# Every row (including the header rows above) whose column L loses its
# "AENDERUNG" text/style and becomes an empty, centered gray cell.
$lOnlyRows = @(13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 38, 39, 40, 41, 42, 43, 71, 99, 103)

$styleDonorRow = 9

foreach ($r in $fullRestyleRows) {
    $src = $ws.Range("A" + $styleDonorRow + ":V" + $styleDonorRow)
    $src.Copy()
    $dst = $ws.Range("A" + $r + ":V" + $r)
    $dst.PasteSpecial($xlPasteFormats)
}

foreach ($r in $lOnlyRows) {
    $srcL = $ws.Range("L" + $styleDonorRow)
    $srcL.Copy()
    $dstL = $ws.Range("L" + $r)
    $dstL.PasteSpecial($xlPasteFormats)
    $dstL.ClearContents()
}

$excel.CutCopyMode = $false
